# ticket583 / #704 #708: "Info" sheet gained a leading "Metadata" label
# column, pushing the existing @base/@data-namespace/... rows one column
# to the right (A:C -> B:D). Re-create that with a real column insert so
# every formula/format/width shifts the way Excel would shift it, then
# stamp the new label into the freed-up column A and park the selection
# back on A1 (it had drifted to C2 before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Insert a new blank column before column A; B:D inherit what used to be
# A:C (values, shared-string refs, per-cell styles and column widths all
# shift automatically, same as picking "Insert" on the column header).
$ws.Columns("A:A").Insert()

# New first column holds the sheet's metadata caption.
$ws.Range("A1").Value = "Metadata"

# Restore the active selection to the top-left cell.
$null = $ws.Range("A1").Select()
